# Auto-generated edit script: updates currentAveragePrice / Leve price / profit
# columns (H-N) for specific leve rows across all 8 crafting-job sheets, per the
# scheduled market-data refresh run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 180.5
$ws.Range("I6").Value = 180.5
$ws.Range("K6").Value = 541.5
$ws.Range("M6").Value = -429.5
$ws.Range("H12").Value = 497.85715
$ws.Range("J12").Value = 94
$ws.Range("L12").Value = 94
$ws.Range("N12").Value = -434
$ws.Range("H33").Value = 252.19048
$ws.Range("I33").Value = 151.52942
$ws.Range("J33").Value = 680
$ws.Range("K33").Value = 151.52942
$ws.Range("L33").Value = 680
$ws.Range("M33").Value = 77.47058000000001
$ws.Range("N33").Value = -1138
$ws.Range("H44").Value = 63888.5
$ws.Range("J44").Value = 63888.5
$ws.Range("L44").Value = 63888.5
$ws.Range("N44").Value = -64812.5
$ws.Range("H62").Value = 149561
$ws.Range("I62").Value = 161874.42
$ws.Range("J62").Value = 1800
$ws.Range("K62").Value = 161874.42
$ws.Range("L62").Value = 1800
$ws.Range("M62").Value = -161250.42
$ws.Range("N62").Value = -3048
$ws.Range("H65").Value = 149561
$ws.Range("I65").Value = 161874.42
$ws.Range("J65").Value = 1800
$ws.Range("K65").Value = 809372.1000000001
$ws.Range("L65").Value = 9000
$ws.Range("M65").Value = -806252.1000000001
$ws.Range("N65").Value = -15240
$ws.Range("H80").Value = 1815.625
$ws.Range("I80").Value = 942.8333
$ws.Range("J80").Value = 2937.7856
$ws.Range("K80").Value = 2828.4999
$ws.Range("L80").Value = 8813.356800000001
$ws.Range("M80").Value = -1830.4999
$ws.Range("N80").Value = -10809.3568
$ws.Range("H83").Value = 1815.625
$ws.Range("I83").Value = 942.8333
$ws.Range("J83").Value = 2937.7856
$ws.Range("K83").Value = 8485.4997
$ws.Range("L83").Value = 26440.0704
$ws.Range("M83").Value = -3493.4997
$ws.Range("N83").Value = -36424.0704
$ws.Range("H98").Value = 1641.3695
$ws.Range("I98").Value = 1685.2972
$ws.Range("J98").Value = 1460.7778
$ws.Range("K98").Value = 1685.2972
$ws.Range("L98").Value = 1460.7778
$ws.Range("M98").Value = -187.2972
$ws.Range("N98").Value = -4456.7778
$ws.Range("H121").Value = 3443
$ws.Range("J121").Value = 3443
$ws.Range("L121").Value = 10329
$ws.Range("N121").Value = -13823
$ws.Range("H122").Value = 1641.3695
$ws.Range("I122").Value = 1685.2972
$ws.Range("J122").Value = 1460.7778
$ws.Range("K122").Value = 5055.8916
$ws.Range("L122").Value = 4382.3334
$ws.Range("M122").Value = -2605.8916
$ws.Range("N122").Value = -9282.3334
$ws.Range("H123").Value = 87000
$ws.Range("I123").Value = 0
$ws.Range("K123").Value = 0
$ws.Range("M123").Value = $null
$ws.Range("H134").Value = 128500
$ws.Range("J134").Value = 128500
$ws.Range("L134").Value = 128500
$ws.Range("N134").Value = -138640
$ws.Range("H137").Value = 2969.6924
$ws.Range("I137").Value = 2934
$ws.Range("J137").Value = 2980.4
$ws.Range("K137").Value = 8802
$ws.Range("L137").Value = 8941.200000000001
$ws.Range("M137").Value = -6252
$ws.Range("N137").Value = -14041.2
$ws.Range("H138").Value = 1253541.9
$ws.Range("I138").Value = 2323.8
$ws.Range("J138").Value = 3338905.2
$ws.Range("K138").Value = 6971.400000000001
$ws.Range("L138").Value = 10016715.6
$ws.Range("M138").Value = -1831.400000000001
$ws.Range("N138").Value = -10026995.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2903.59
$ws.Range("I32").Value = 2882.8447
$ws.Range("K32").Value = 2882.8447
$ws.Range("M32").Value = -2595.8447
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").Value = $null
$ws.Range("H45").Value = 3787
$ws.Range("I45").Value = 2938.4285
$ws.Range("K45").Value = 2938.4285
$ws.Range("M45").Value = -2561.4285
$ws.Range("H97").Value = 5514.615
$ws.Range("I97").Value = 2586.85
$ws.Range("K97").Value = 2586.85
$ws.Range("M97").Value = -2090.85

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H24").Value = 655
$ws.Range("I24").Value = 655
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 655
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = -420
$ws.Range("N24").Value = $null
$ws.Range("H86").Value = 1808940.9
$ws.Range("I86").Value = 2675147.8
$ws.Range("J86").Value = 4343.25
$ws.Range("K86").Value = 2675147.8
$ws.Range("L86").Value = 4343.25
$ws.Range("M86").Value = -2674024.8
$ws.Range("N86").Value = -6589.25
$ws.Range("H89").Value = 1808940.9
$ws.Range("I89").Value = 2675147.8
$ws.Range("J89").Value = 4343.25
$ws.Range("K89").Value = 13375739
$ws.Range("L89").Value = 21716.25
$ws.Range("M89").Value = -13370123
$ws.Range("N89").Value = -32948.25
$ws.Range("H94").Value = 1451.174
$ws.Range("I94").Value = 1550
$ws.Range("K94").Value = 1550
$ws.Range("M94").Value = -1099
$ws.Range("H99").Value = 2108.7827
$ws.Range("I99").Value = 2066.8572
$ws.Range("J99").Value = 2549
$ws.Range("K99").Value = 2066.8572
$ws.Range("L99").Value = 2549
$ws.Range("M99").Value = -568.8571999999999
$ws.Range("N99").Value = -5545

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 183999.33
$ws.Range("J9").Value = 183999.33
$ws.Range("L9").Value = 183999.33
$ws.Range("N9").Value = -184335.33
$ws.Range("H19").Value = 1812.1
$ws.Range("I19").Value = 202.625
$ws.Range("J19").Value = 8250
$ws.Range("K19").Value = 202.625
$ws.Range("L19").Value = 8250
$ws.Range("M19").Value = -32.625
$ws.Range("N19").Value = -8590
$ws.Range("H24").Value = 1812.1
$ws.Range("I24").Value = 202.625
$ws.Range("J24").Value = 8250
$ws.Range("K24").Value = 202.625
$ws.Range("L24").Value = 8250
$ws.Range("M24").Value = -32.625
$ws.Range("N24").Value = -8590
$ws.Range("H31").Value = 5501.5
$ws.Range("I31").Value = 3605.7646
$ws.Range("J31").Value = 6541.0967
$ws.Range("K31").Value = 3605.7646
$ws.Range("L31").Value = 6541.0967
$ws.Range("M31").Value = -3310.7646
$ws.Range("N31").Value = -7131.0967
$ws.Range("H34").Value = 5501.5
$ws.Range("I34").Value = 3605.7646
$ws.Range("J34").Value = 6541.0967
$ws.Range("K34").Value = 3605.7646
$ws.Range("L34").Value = 6541.0967
$ws.Range("M34").Value = -3403.7646
$ws.Range("N34").Value = -6945.0967
$ws.Range("H86").Value = 8329.809999999999
$ws.Range("I86").Value = 5475.5454
$ws.Range("J86").Value = 11469.5
$ws.Range("K86").Value = 5475.5454
$ws.Range("L86").Value = 11469.5
$ws.Range("M86").Value = -4352.5454
$ws.Range("N86").Value = -13715.5
$ws.Range("H89").Value = 8329.809999999999
$ws.Range("I89").Value = 5475.5454
$ws.Range("J89").Value = 11469.5
$ws.Range("K89").Value = 27377.727
$ws.Range("L89").Value = 57347.5
$ws.Range("M89").Value = -21761.727
$ws.Range("N89").Value = -68579.5
$ws.Range("H132").Value = 1972.8422
$ws.Range("I132").Value = 1992.3077
$ws.Range("J132").Value = 1930.6666
$ws.Range("K132").Value = 5976.9231
$ws.Range("L132").Value = 5791.9998
$ws.Range("M132").Value = -3446.9231
$ws.Range("N132").Value = -10851.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1638.238
$ws.Range("I34").Value = 1570
$ws.Range("K34").Value = 4710
$ws.Range("M34").Value = -4626
$ws.Range("H56").Value = 1304557.5
$ws.Range("I56").Value = 1304557.5
$ws.Range("K56").Value = 1304557.5
$ws.Range("M56").Value = -1304027.5
$ws.Range("H68").Value = 7269.857
$ws.Range("I68").Value = 1750
$ws.Range("J68").Value = 9477.799999999999
$ws.Range("K68").Value = 5250
$ws.Range("L68").Value = 28433.4
$ws.Range("M68").Value = -4439
$ws.Range("N68").Value = -30055.4
$ws.Range("H71").Value = 7269.857
$ws.Range("I71").Value = 1750
$ws.Range("J71").Value = 9477.799999999999
$ws.Range("K71").Value = 15750
$ws.Range("L71").Value = 85300.2
$ws.Range("M71").Value = -11694
$ws.Range("N71").Value = -93412.2
$ws.Range("H114").Value = 5185
$ws.Range("I114").Value = 5000
$ws.Range("J114").Value = 5555
$ws.Range("K114").Value = 15000
$ws.Range("L114").Value = 16665
$ws.Range("M114").Value = -11746
$ws.Range("N114").Value = -23173

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 69849.266
$ws.Range("I43").Value = 28006.4
$ws.Range("J43").Value = 90770.7
$ws.Range("K43").Value = 28006.4
$ws.Range("L43").Value = 90770.7
$ws.Range("M43").Value = -27855.4
$ws.Range("N43").Value = -91072.7
$ws.Range("H46").Value = 27545.455
$ws.Range("J46").Value = 30000
$ws.Range("L46").Value = 30000
$ws.Range("N46").Value = -30312
$ws.Range("H57").Value = 90844.60000000001
$ws.Range("J57").Value = 90844.60000000001
$ws.Range("L57").Value = 90844.60000000001
$ws.Range("N57").Value = -92484.60000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5959.5454
$ws.Range("J40").Value = 6536.4287
$ws.Range("L40").Value = 6536.4287
$ws.Range("N40").Value = -6808.4287
$ws.Range("H93").Value = 1220
$ws.Range("I93").Value = 1220
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 1220
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = 28
$ws.Range("N93").Value = $null
$ws.Range("H136").Value = 5348.069
$ws.Range("I136").Value = 5458.8184
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 16376.4552
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -13826.4552
$ws.Range("N136").Value = -20100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 4700
$ws.Range("I96").Value = 4100
$ws.Range("J96").Value = 5000
$ws.Range("K96").Value = 4100
$ws.Range("L96").Value = 5000
$ws.Range("M96").Value = -2727
$ws.Range("N96").Value = -7746
$ws.Range("H122").Value = 4244.2
$ws.Range("I122").Value = 2959.3635
$ws.Range("J122").Value = 13666.333
$ws.Range("K122").Value = 8878.0905
$ws.Range("L122").Value = 40998.999
$ws.Range("M122").Value = -6428.0905
$ws.Range("N122").Value = -45898.999
$ws.Range("H132").Value = 1964.25
$ws.Range("I132").Value = 1752.3077
$ws.Range("J132").Value = 2357.8572
$ws.Range("K132").Value = 5256.9231
$ws.Range("L132").Value = 7073.571599999999
$ws.Range("M132").Value = -2726.9231
$ws.Range("N132").Value = -12133.5716
